$wb = $excel.ActiveWorkbook

# --- Trucks sheet: add a "fuel cost" column (E) ------------------------
$trucks = $wb.Worksheets.Item("Trucks")
$trucks.Range("E1").Value = "fuel cost"
for ($r = 2; $r -le 11; $r++) {
    $trucks.Cells.Item($r, 5).Value = 12
}

# --- Others sheet: add a "labor cost" row (row 6) -----------------------
$others = $wb.Worksheets.Item("Others")
$others.Range("A6").Value = "labor cost"
$others.Range("B6").Value = 25

# --- Make "Others" the active sheet/tab, with A7 selected --------------
$others.Activate() | Out-Null
$others.Range("A7").Select() | Out-Null
